$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "deuteron" target value with "d" for data rows 2-10 (column I = target)
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value2 -eq "deuteron") {
        $cell.Value = "d"
    }
}

# Make header row bold (row 1, columns A1:K1), keep existing center alignment
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true

# Update selection to match target state
$ws.Range("C15").Select()
